$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: A1 stays "Animal_ID" (string index changes internally but displayed text is same)
$ws.Range("A1").Value = "Animal_ID"

# Update Group labels: rows 2-8 were "A" -> now "Group A"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = "Group A"
}

# Update Group labels: rows 9-15 were "B" -> now "Group B"
for ($r = 9; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Value = "Group B"
}

# Animal_ID numeric changes
$ws.Range("A6").Value = 67
$ws.Range("A13").Value = 40

# Temperature data changes in columns C and D for rows 9-15
$ws.Range("C9").Value = 39
$ws.Range("D9").Value = 40

$ws.Range("C10").Value = 37.1
$ws.Range("D10").Value = 39.7

$ws.Range("C11").Value = 38
$ws.Range("D11").Value = 39.5

$ws.Range("C12").Value = 38.5
$ws.Range("D12").Value = 39.5

$ws.Range("C13").Value = 37.8
$ws.Range("D13").Value = 38

$ws.Range("C14").Value = 39
$ws.Range("D14").Value = 39.7

$ws.Range("C15").Value = 38.8
$ws.Range("D15").Value = 39.8

# Update the active selection to match the edited workbook state
$ws.Range("D16").Select()
